$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings usage: column D (Target cluster) values change, and the
# row count shrinks from 6 data rows (rows 2-7) down to 3 data rows (rows 2-4)
# with freshly recomputed TPM-based statistics.

# Row 2: MuSCs -> Avp -> Avpr1a -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Avp"
$ws.Range("C2").Value = "Avpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.2600935
$ws.Range("H2").Value = 0.520187
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.5383755
$ws.Range("N2").Value = 5.076751
$ws.Range("O2").Value = 0.1536689326353429
$ws.Range("P2").Value = 0.126906281076628
$ws.Range("Q2").Value = 0.66021496810925
$ws.Range("R2").Value = 2.640859872437
$ws.Range("S2").Value = 0.1536689326353429
$ws.Range("T2").Value = 0.126906281076628

# Row 3: MuSCs -> Avp -> Avpr1a -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Avp"
$ws.Range("C3").Value = "Avpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.2600935
$ws.Range("H3").Value = 0.520187
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.967000000000001
$ws.Range("N3").Value = 20.901
$ws.Range("O3").Value = 0.4217703226612587
$ws.Range("P3").Value = 0.5224735624777741
$ws.Range("Q3").Value = 1.8120714145
$ws.Range("R3").Value = 10.872428487
$ws.Range("S3").Value = 0.4217703226612587
$ws.Range("T3").Value = 0.5224735624777741

# Row 4: MuSCs -> Avp -> Avpr1a -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Avp"
$ws.Range("C4").Value = "Avpr1a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.2600935
$ws.Range("H4").Value = 0.520187
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.0130935
$ws.Range("N4").Value = 14.026187
$ws.Range("O4").Value = 0.4245607447033983
$ws.Range("P4").Value = 0.3506201564455979
$ws.Range("Q4").Value = 1.82406003424225
$ws.Range("R4").Value = 7.296240136969
$ws.Range("S4").Value = 0.4245607447033983
$ws.Range("T4").Value = 0.3506201564455979

# Remove rows 5-7, which held the old "Neutrophils" sending-cluster data that
# no longer exists after the TPM recompute.
$ws.Range("A5:T7").Delete()
